{"js": "// \"revisi penulisan teorema + bab 5 uji coba\"\n// The visible/textual change in this document is limited to the\n// \"Lembar Pengesahan\" (approval page) examiner table: the three\n// placeholder \"________\" examiner-name blanks are filled in with the\n// actual examiners' names.\n//\n//   Penguji I   : \"________\"                -> \"Dr. Yosi Kristian\"\n//   Penguji II  : \"________, M.App.Sc\"      -> \"Dr. Ir. Esther Irawati Setiawan, S.Kom, M.Kom\"\n//   Penguji III : \"________\"                -> \"Evan Kusuma Susanto\"\n//\n// Each paragraph ends with \"(Penguji I)\" / \"(Penguji II)\" / \"(Penguji III)\"\n// which we use to unambiguously locate the right row before doing a\n// scoped search-and-replace inside that paragraph only (so the three\n// identical \"________\" placeholders don't collide with one another).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nasync function replaceInParagraph(paragraph, findText, replaceText) {\n  const results = paragraph.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text;\n\n  if (text.indexOf(\"(Penguji I)\") !== -1) {\n    await replaceInParagraph(paragraph, \"________\", \"Dr. Yosi Kristian\");\n  } else if (text.indexOf(\"(Penguji II)\") !== -1) {\n    await replaceInParagraph(\n      paragraph,\n      \"________, M.App.Sc\",\n      \"Dr. Ir. Esther Irawati Setiawan, S.Kom, M.Kom\"\n    );\n  } else if (text.indexOf(\"(Penguji III)\") !== -1) {\n    await replaceInParagraph(paragraph, \"________\", \"Evan Kusuma Susanto\");\n  }\n}\n", "ps1": "# \"revisi penulisan teorema + bab 5 uji coba\"\n# The visible/textual change in this document is limited to the\n# \"Lembar Pengesahan\" (approval page) examiner table: the three\n# placeholder \"________\" examiner-name blanks are filled in with the\n# actual examiners' names.\n#\n#   Penguji I   : \"________\"            -> \"Dr. Yosi Kristian\"\n#   Penguji II  : \"________, M.App.Sc\"  -> \"Dr. Ir. Esther Irawati Setiawan, S.Kom, M.Kom\"\n#   Penguji III : \"________\"            -> \"Evan Kusuma Susanto\"\n#\n# Each paragraph ends with \"(Penguji I)\" / \"(Penguji II)\" / \"(Penguji III)\"\n# which we use to unambiguously locate the right row, then scope the\n# Find/Replace (wdReplaceOne) to that paragraph's own Range so the three\n# identical \"________\" placeholders don't collide with one another.\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n\n    if ($text -like \"*(Penguji I)*\") {\n        $rng = $p.Range\n        $rng.Find.Execute(\"________\", $false, $false, $false, $false, $false, $true, 1, $false, \"Dr. Yosi Kristian\", 1) | Out-Null\n    }\n    elseif ($text -like \"*(Penguji II)*\") {\n        $rng = $p.Range\n        $rng.Find.Execute(\"________, M.App.Sc\", $false, $false, $false, $false, $false, $true, 1, $false, \"Dr. Ir. Esther Irawati Setiawan, S.Kom, M.Kom\", 1) | Out-Null\n    }\n    elseif ($text -like \"*(Penguji III)*\") {\n        $rng = $p.Range\n        $rng.Find.Execute(\"________\", $false, $false, $false, $false, $false, $true, 1, $false, \"Evan Kusuma Susanto\", 1) | Out-Null\n    }\n}\n"}
